# Backend import update for simplified diagnostic
# Renames/re-orders the "quality" columns on the header row and appends
# four new "EGAlim" detail columns (viandes/produits aquatiques).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) --------------------------------------------------
# Columns A:N are unchanged in content. O:R are renamed in place, and four
# brand-new columns S:V are appended.
$ws.Range("O1").Value = "SIQO"
$ws.Range("P1").Value = "Environnement"
$ws.Range("Q1").Value = "Autre EGAlim"
$ws.Range("R1").Value = "Viandes vollailes total"
$ws.Range("S1").Value = "Viandes vollailles EGAlim"
$ws.Range("T1").Value = "Viandes vollailles provenance France"
$ws.Range("U1").Value = "Produits aquatiques total"
$ws.Range("V1").Value = "Produits aquatiques EGAlim"

# The whole header row is now bold.
$ws.Range("A1:V1").Font.Bold = $true

# --- Column widths ---------------------------------------------------------
# Column O shrinks to the same width as the "Total"/"Bio" columns, and the
# newly inserted column P gets a narrower, auto-fit-ish width.
$ws.Columns.Item(15).ColumnWidth = 4.5
$ws.Columns.Item(16).ColumnWidth = 12.83

# --- Selection --------------------------------------------------------------
$null = $ws.Range("V1").Select()
